$wb = $excel.ActiveWorkbook

# New BSE-listed names that showed up in the latest stock.yaml "week"
# scrape (run at 07/02/2025 11:32:53) — sr, nsecode, name, bsecode,
# per_chg, close, volume, timeframe, Date Time.
$data = @(
    ,@(1, "OFSS", "Oracle Financial Services Software Limited", "532466", 1.93, 9452.700000000001, 111132, "week", "07/02/2025 11:32:53")
    ,@(2, "COFORGE", "Coforge (Niit Tech)", "532541", 0.17, 8535.85, 301753, "week", "07/02/2025 11:32:53")
    ,@(3, "LTIM", "LTI Mindtree Ltd", "540005", 0.2, 5956.6, 96126, "week", "07/02/2025 11:32:53")
    ,@(4, "TRENT", "Trent Limited", "500251", 3.36, 5454.4, 2846965, "week", "07/02/2025 11:32:53")
    ,@(5, "NAVINFLUOR", "Navin Fluorine International Limited", "532504", -0.65, 4220, 209989, "week", "07/02/2025 11:32:53")
    ,@(6, "DMART", "Avenue Supermarts", "540376", -0.13, 3743.5, 301739, "week", "07/02/2025 11:32:53")
    ,@(7, "BALKRISIND", "Balkrishna Industries Limited", "502355", -1.44, 2722.9, 102316, "week", "07/02/2025 11:32:53")
    ,@(8, "ADANIENT", "Adani Enterprises Limited", "512599", 1.14, 2338.8, 892900, "week", "07/02/2025 11:32:53")
    ,@(9, "LUPIN", "Lupin Limited", "500257", 0.72, 2201.15, 776517, "week", "07/02/2025 11:32:53")
    ,@(10, "ICICIGI", "ICICI Lombard General Insurance Company Ltd", "540716", -0.1, 1826, 352127, "week", "07/02/2025 11:32:53")
    ,@(11, "SUNPHARMA", "Sun Pharmaceuticals Industries Limited", "524715", 0.45, 1749.8, 1371864, "week", "07/02/2025 11:32:53")
    ,@(12, "POLICYBZR", "PB Fintech Ltd", "543390", 0.95, 1718, 916722, "week", "07/02/2025 11:32:53")
    ,@(13, "UNITDSPR", "United Spirits Ltd", "532432", 0.88, 1427.5, 621363, "week", "07/02/2025 11:32:53")
    ,@(14, "VOLTAS", "Voltas Limited", "500575", 2.33, 1383.25, 1597612, "week", "07/02/2025 11:32:53")
    ,@(15, "AUROPHARMA", "Aurobindo Pharma Limited", "524804", -0.01, 1191.7, 1233173, "week", "07/02/2025 11:32:53")
    ,@(16, "GODREJCP", "Godrej Consumer Products Limited", "532424", -1.15, 1106.75, 478213, "week", "07/02/2025 11:32:53")
    ,@(17, "SBICARD", "SBI Cards & Payment Services Ltd", "543066", 0.6, 815.55, 1378212, "week", "07/02/2025 11:32:53")
    ,@(18, "INDHOTEL", "The Indian Hotels Company Limited", "500850", 0.22, 791.8, 1439580, "week", "07/02/2025 11:32:53")
    ,@(19, "UPL", "Upl Limited", "512070", -0.53, 641.35, 2340663, "week", "07/02/2025 11:32:53")
    ,@(20, "VBL", "Varun Beverages Limited", "540180", -2.62, 553.7, 6361141, "week", "07/02/2025 11:32:53")
    ,@(21, "KALYANKJIL", "Kalyan Jewellers India Ltd", "543278", -0.55, 539.4, 17368688, "week", "07/02/2025 11:32:53")
    ,@(22, "BERGEPAINT", "Berger Paints (i) Limited", "509480", 0.7, 479.35, 373238, "week", "07/02/2025 11:32:53")
    ,@(23, "ITC", "Itc Limited", "500875", -2.32, 430.85, 18875982, "week", "07/02/2025 11:32:53")
    ,@(24, "PETRONET", "Petronet Lng Limited", "532522", 1.59, 315.45, 2157166, "week", "07/02/2025 11:32:53")
    ,@(25, "ONGC", "Oil & Natural Gas Corporation Limited", "500312", -2.83, 248.9, 10098454, "week", "07/02/2025 11:32:53")
    ,@(26, "ASHOKLEY", "Ashok Leyland Limited", "500477", 1.08, 212.05, 3939048, "week", "07/02/2025 11:32:53")
    ,@(27, "IGL", "Indraprastha Gas Limited", "532514", -1.7, 198.9, 1635563, "week", "07/02/2025 11:32:53")
)

# --- Part 1: "day" sheet — D1547:D1555 were numeric BSE codes stored as
# text (inlineStr); convert them to genuine numbers. ---
$wsDay = $wb.Worksheets.Item("day")
$wsDay.Range("D1547").Value = 533179
$wsDay.Range("D1548").Value = 526299
$wsDay.Range("D1549").Value = 524715
$wsDay.Range("D1550").Value = 532755
$wsDay.Range("D1551").Value = 500493
$wsDay.Range("D1552").Value = 532868
$wsDay.Range("D1553").Value = 500425
$wsDay.Range("D1554").Value = 539336
$wsDay.Range("D1555").Value = 524000

# --- Part 2: "week" sheet — append 27 new rows (923-949) pulled in by the
# latest stock.yaml refresh (07/02/2025 run). bsecode (column D) keeps the
# text/inlineStr representation used by the scraper for this sheet. ---
$wsWeek = $wb.Worksheets.Item("week")
$startRow = 923
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $startRow + $idx
    $rec = $data[$idx]
    $wsWeek.Cells.Item($row, 1).Value = $rec[0]
    $wsWeek.Cells.Item($row, 2).Value = $rec[1]
    $wsWeek.Cells.Item($row, 3).Value = $rec[2]
    $wsWeek.Cells.Item($row, 4).NumberFormat = "@"
    $wsWeek.Cells.Item($row, 4).Value = $rec[3]
    $wsWeek.Cells.Item($row, 4).Style = "Normal"
    $wsWeek.Cells.Item($row, 5).Value = $rec[4]
    $wsWeek.Cells.Item($row, 6).Value = $rec[5]
    $wsWeek.Cells.Item($row, 7).Value = $rec[6]
    $wsWeek.Cells.Item($row, 8).Value = $rec[7]
    $wsWeek.Cells.Item($row, 9).Value = $rec[8]
}
